# feat: add 2022-Q3 data
#
# The "2021-Q4" sheet's data becomes the new "2022-Q3" data, and the old
# "2021-Q4" data is preserved on a fresh sheet of its own (inserted right
# after, so it keeps the sheet in the same relative spot). The "总计"
# (totals) sheet gets a new leading row for 2022-Q3 and keeps the old
# 2021-Q4 row underneath it.

function Set-TextCell($cell, $value) {
    # Plain `.Value = "..."` lets Excel infer the type, so numeric-looking
    # strings (fund codes like "167703", percentages like "0.30") would be
    # silently coerced to numbers and lose their formatting (leading
    # zeros / trailing zeros). Forcing a text number format first keeps the
    # literal string; resetting the style afterwards avoids leaving the
    # cell permanently tagged with a "text" style it didn't have before.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

$totals = $wb.Worksheets.Item(1)
$q4sheet = $wb.Worksheets.Item(2)

# --- Duplicate the current "2021-Q4" sheet so its existing data is kept on
# --- its own sheet, positioned right after the (soon to be renamed) sheet.
$q4sheet.Copy($null, $q4sheet)
$oldDataSheet = $wb.Worksheets.Item(3)

# The original sheet becomes the new quarter's sheet; the copy keeps the
# old quarter's name and data untouched.
$q4sheet.Name = "2022-Q3"
$oldDataSheet.Name = "2021-Q4"

# --- Overwrite the (renamed) sheet with the new 2022-Q3 fund data.
$newSheet = $q4sheet

Set-TextCell $newSheet.Cells.Item(2, 2) "167703"
Set-TextCell $newSheet.Cells.Item(2, 3) "德邦量化优选股票（LOF）C"
Set-TextCell $newSheet.Cells.Item(2, 4) "0.30"
Set-TextCell $newSheet.Cells.Item(2, 5) "90.17"
Set-TextCell $newSheet.Cells.Item(2, 6) "0.93"
Set-TextCell $newSheet.Cells.Item(2, 7) "0.0028"
$newSheet.Cells.Item(2, 8).Value = 9

Set-TextCell $newSheet.Cells.Item(3, 2) "167702"
Set-TextCell $newSheet.Cells.Item(3, 3) "德邦量化优选股票（LOF）A"
Set-TextCell $newSheet.Cells.Item(3, 4) "0.17"
Set-TextCell $newSheet.Cells.Item(3, 5) "90.17"
Set-TextCell $newSheet.Cells.Item(3, 6) "0.93"
Set-TextCell $newSheet.Cells.Item(3, 7) "0.0016"
$newSheet.Cells.Item(3, 8).Value = 9

# --- Update the "总计" summary sheet: push the old 2021-Q4 row down to row
# --- 3, and turn row 2 into the new 2022-Q3 entry.
# Row 3 is brand new, so A3 needs to pick up the same "index column" style
# A2 already has (bold, centered, bordered) - copy the formatting across
# before writing the value.
$totals.Cells.Item(2, 1).Copy($totals.Cells.Item(3, 1))
$totals.Cells.Item(3, 1).Value = 1
Set-TextCell $totals.Cells.Item(3, 2) "2021-Q4"
$totals.Cells.Item(3, 3).Value = 2
$totals.Cells.Item(3, 4).Value = 0.02

Set-TextCell $totals.Cells.Item(2, 2) "2022-Q3"
$totals.Cells.Item(2, 3).Value = 2
$totals.Cells.Item(2, 4).Value = 0

# The sheet-copy above shifts Excel's active-sheet selection; restore it to
# the "总计" sheet so the workbook-level view state is left as it was found.
$totals.Activate()
